# Lower score for Tiny Toons 2
# Re-rate a few achievements' difficulty (lowering their point value) and
# mark the Checklist rows that are now fully done.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Checklist sheet: tick off column H ("Triggers?") for the rows that are
# now fully done, which flips their column I formula result to "YES".
# Do this sheet first so its own window selection is recorded before we
# flip back to Achievements (the sheet that should stay active).
# ----------------------------------------------------------------------
$wsChk = $wb.Worksheets.Item("Checklist")
$wsChk.Activate()

$rows = @(12,13,14,15,16,25,26,27,28,29,30,31,33)
foreach ($r in $rows) {
    $wsChk.Cells.Item($r, 8).Value = "X"
}

# restore the selection shown in the file for this sheet
$wsChk.Range("H11:H16").Select()

# ----------------------------------------------------------------------
# Achievements sheet: re-rate a few achievements' difficulty.
# Column D holds the difficulty text; column E is a VLOOKUP formula that
# derives the point value from the Stats sheet, so it recalculates on
# its own once D changes. All the Stats/Extras/Text sheet totals are
# themselves formulas driven off this data, so they recalc too.
# ----------------------------------------------------------------------
$wsAch = $wb.Worksheets.Item("Achievements")
$wsAch.Activate()

$wsAch.Range("D2").Value = "Hard"       # was Very Easy (2 pts) -> Hard (5 pts)
$wsAch.Range("D3").Value = "Hard"       # was Easy (3 pts) -> Hard (5 pts)
$wsAch.Range("D27").Value = "Very Hard" # was Super Hard (25 pts) -> Very Hard (10 pts)
$wsAch.Range("D28").Value = "Very Hard" # was Super Hard (25 pts) -> Very Hard (10 pts)

# restore the selection shown in the file for this sheet (and keep it the
# active sheet/tab, matching the original workbook)
$wsAch.Range("E2:E33").Select()
